$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing data row (151) onto each new row,
# then set the values for the cells that the source data provides.

# Row 152 (movie #151) - Sexta feira 13 parte 5
$ws.Range("A151:J151").Copy()
$ws.Range("A152:J152").PasteSpecial(-4122)
$ws.Range("A152").Value = 151
$ws.Range("B152").Value = "Sexta feira 13 parte 5"
$ws.Range("C152").Value = "7,5"
$ws.Range("D152").Value = "8,5"

# Row 153 (movie #152) - Sexta feira 13 parte 6
$ws.Range("A151:J151").Copy()
$ws.Range("A153:J153").PasteSpecial(-4122)
$ws.Range("A153").Value = 152
$ws.Range("B153").Value = "Sexta feira 13 parte 6"
$ws.Range("C153").Value = "7"
$ws.Range("D153").Value = "7"

# Row 154 (movie #153) - Sexta feira 13 parte 7
$ws.Range("A151:J151").Copy()
$ws.Range("A154:J154").PasteSpecial(-4122)
$ws.Range("A154").Value = 153
$ws.Range("B154").Value = "Sexta feira 13 parte 7"
$ws.Range("C154").Value = "7,5"
$ws.Range("D154").Value = "8,5"

# Row 155 (movie #154) - Sexta feira 13 parte 8
$ws.Range("A151:J151").Copy()
$ws.Range("A155:J155").PasteSpecial(-4122)
$ws.Range("A155").Value = 154
$ws.Range("B155").Value = "Sexta feira 13 parte 8"
$ws.Range("C155").Value = "6"
$ws.Range("D155").Value = "8"

# Row 156 (movie #155) - Sexta feira 13 parte 9
$ws.Range("A151:J151").Copy()
$ws.Range("A156:J156").PasteSpecial(-4122)
$ws.Range("A156").Value = 155
$ws.Range("B156").Value = "Sexta feira 13 parte 9"
$ws.Range("C156").Value = "6"
$ws.Range("D156").Value = "7"

# Row 157 (movie #156) - Sexta feira 13 parte 10
$ws.Range("A151:J151").Copy()
$ws.Range("A157:J157").PasteSpecial(-4122)
$ws.Range("A157").Value = 156
$ws.Range("B157").Value = "Sexta feira 13 parte 10"
$ws.Range("C157").Value = "6"
$ws.Range("D157").Value = "10"

# Row 158 (movie #157) - Sexta feira 13 parte 11
$ws.Range("A151:J151").Copy()
$ws.Range("A158:J158").PasteSpecial(-4122)
$ws.Range("A158").Value = 157
$ws.Range("B158").Value = "Sexta feira 13 parte 11"
$ws.Range("C158").Value = "9"
$ws.Range("D158").Value = "8,5"

# Row 159 (movie #158) - Sexta feira 13 remake
$ws.Range("A151:J151").Copy()
$ws.Range("A159:J159").PasteSpecial(-4122)
$ws.Range("A159").Value = 158
$ws.Range("B159").Value = "Sexta feira 13 remake"
$ws.Range("C159").Value = "8"
$ws.Range("D159").Value = "8,5"

# Row 160 (movie #159) - Halloween
$ws.Range("A151:J151").Copy()
$ws.Range("A160:J160").PasteSpecial(-4122)
$ws.Range("A160").Value = 159
$ws.Range("B160").Value = "Halloween"
$ws.Range("C160").Value = "8"
$ws.Range("D160").Value = "10"

# Row 161 (movie #160) - Halloween 2
$ws.Range("A151:J151").Copy()
$ws.Range("A161:J161").PasteSpecial(-4122)
$ws.Range("A161").Value = 160
$ws.Range("B161").Value = "Halloween 2"
$ws.Range("C161").Value = "8"
$ws.Range("D161").Value = "8,5"

# Row 162 (movie #161) - Halloween 3
$ws.Range("A151:J151").Copy()
$ws.Range("A162:J162").PasteSpecial(-4122)
$ws.Range("A162").Value = 161
$ws.Range("B162").Value = "Halloween 3"
$ws.Range("C162").Value = "6,5"
$ws.Range("D162").Value = "6"

# Row 163 (movie #162) - Halloween 4
$ws.Range("A151:J151").Copy()
$ws.Range("A163:J163").PasteSpecial(-4122)
$ws.Range("A163").Value = 162
$ws.Range("B163").Value = "Halloween 4"
$ws.Range("C163").Value = "9"
$ws.Range("D163").Value = "9"

# Row 164 (movie #163) - Halloween 5
$ws.Range("A151:J151").Copy()
$ws.Range("A164:J164").PasteSpecial(-4122)
$ws.Range("A164").Value = 163
$ws.Range("B164").Value = "Halloween 5"
$ws.Range("C164").Value = "8"
$ws.Range("D164").Value = "6,5"

# Row 165 (movie #164) - Halloween 6
$ws.Range("A151:J151").Copy()
$ws.Range("A165:J165").PasteSpecial(-4122)
$ws.Range("A165").Value = 164
$ws.Range("B165").Value = "Halloween 6"
$ws.Range("C165").Value = "7"
$ws.Range("D165").Value = "6,5"

# Row 166 (movie #165) - Halloween remake 1
$ws.Range("A151:J151").Copy()
$ws.Range("A166:J166").PasteSpecial(-4122)
$ws.Range("A166").Value = 165
$ws.Range("B166").Value = "Halloween remake 1"
$ws.Range("C166").Value = "7"
$ws.Range("D166").Value = "7,5"
$ws.Range("F166").Value = "5"

# Row 167 (movie #166) - Donnie darko
$ws.Range("A151:J151").Copy()
$ws.Range("A167:J167").PasteSpecial(-4122)
$ws.Range("A167").Value = 166
$ws.Range("B167").Value = "Donnie darko"
$ws.Range("C167").Value = "10"
$ws.Range("D167").Value = "8,5"
$ws.Range("F167").Value = "9"

# Row 168 (movie #167) - Halloween remake 2
$ws.Range("A151:J151").Copy()
$ws.Range("A168:J168").PasteSpecial(-4122)
$ws.Range("A168").Value = 167
$ws.Range("B168").Value = "Halloween remake 2"
$ws.Range("C168").Value = "5,5"
$ws.Range("D168").Value = "4"
$ws.Range("F168").Value = "2"

# Row 169 (movie #168) - Halloween remake (2018)
$ws.Range("A151:J151").Copy()
$ws.Range("A169:J169").PasteSpecial(-4122)
$ws.Range("A169").Value = 168
$ws.Range("B169").Value = "Halloween remake (2018)"
$ws.Range("C169").Value = "9"
$ws.Range("D169").Value = "9,5"

# Row 170 (movie #169) - predador
$ws.Range("A151:J151").Copy()
$ws.Range("A170:J170").PasteSpecial(-4122)
$ws.Range("A170").Value = 169
$ws.Range("B170").Value = "predador"
$ws.Range("C170").Value = "6"
$ws.Range("D170").Value = "8"
$ws.Range("F170").Value = "6"
$ws.Range("G170").Value = "5"
$ws.Range("H170").Value = "6,5"
$ws.Range("I170").Value = "9"

# Row 171 (movie #170) - predador 2
$ws.Range("A151:J151").Copy()
$ws.Range("A171:J171").PasteSpecial(-4122)
$ws.Range("A171").Value = 170
$ws.Range("B171").Value = "predador 2"
$ws.Range("C171").Value = "7"
$ws.Range("D171").Value = "8"

# Row 172 (movie #171) - alien
$ws.Range("A151:J151").Copy()
$ws.Range("A172:J172").PasteSpecial(-4122)
$ws.Range("A172").Value = 171
$ws.Range("B172").Value = "alien"
$ws.Range("C172").Value = "7"
$ws.Range("D172").Value = "7,5"
$ws.Range("G172").Value = "6"

# Row 173 (movie #172) - alien 2
$ws.Range("A151:J151").Copy()
$ws.Range("A173:J173").PasteSpecial(-4122)
$ws.Range("A173").Value = 172
$ws.Range("B173").Value = "alien 2"
$ws.Range("C173").Value = "7,5"
$ws.Range("D173").Value = "7,5"

# Row 174 (movie #173) - alien 3
$ws.Range("A151:J151").Copy()
$ws.Range("A174:J174").PasteSpecial(-4122)
$ws.Range("A174").Value = 173
$ws.Range("B174").Value = "alien 3"
$ws.Range("C174").Value = "8,5"
$ws.Range("D174").Value = "8,5"

# Row 175 (movie #174) - alien 4
$ws.Range("A151:J151").Copy()
$ws.Range("A175:J175").PasteSpecial(-4122)
$ws.Range("A175").Value = 174
$ws.Range("B175").Value = "alien 4"
$ws.Range("C175").Value = "5"
$ws.Range("D175").Value = "3,5"

# Row 176 (movie #175) - alien vs predador
$ws.Range("A151:J151").Copy()
$ws.Range("A176:J176").PasteSpecial(-4122)
$ws.Range("A176").Value = 175
$ws.Range("B176").Value = "alien vs predador"
$ws.Range("C176").Value = "8"
$ws.Range("D176").Value = "7"

# Row 177 (movie #176) - alien vs predador 2
$ws.Range("A151:J151").Copy()
$ws.Range("A177:J177").PasteSpecial(-4122)
$ws.Range("A177").Value = 176
$ws.Range("B177").Value = "alien vs predador 2"
$ws.Range("C177").Value = "6"
$ws.Range("D177").Value = "6"

# Row 178 (movie #177) - a morte de stalin
$ws.Range("A151:J151").Copy()
$ws.Range("A178:J178").PasteSpecial(-4122)
$ws.Range("A178").Value = 177
$ws.Range("B178").Value = "a morte de stalin"
$ws.Range("C178").Value = "9"
$ws.Range("D178").Value = "8,5"
$ws.Range("F178").Value = "9"

# Row 179 (movie #178) - nosferatu
$ws.Range("A151:J151").Copy()
$ws.Range("A179:J179").PasteSpecial(-4122)
$ws.Range("A179").Value = 178
$ws.Range("B179").Value = "nosferatu"
$ws.Range("C179").Value = "6"
$ws.Range("D179").Value = "5"
$ws.Range("F179").Value = "6"

# Row 180 (movie #179) - the trip
$ws.Range("A151:J151").Copy()
$ws.Range("A180:J180").PasteSpecial(-4122)
$ws.Range("A180").Value = 179
$ws.Range("B180").Value = "the trip"
$ws.Range("C180").Value = "8"
$ws.Range("D180").Value = "6"

# Row 181 (movie #180) - encontro de casais
$ws.Range("A151:J151").Copy()
$ws.Range("A181:J181").PasteSpecial(-4122)
$ws.Range("A181").Value = 180
$ws.Range("B181").Value = "encontro de casais"
$ws.Range("C181").Value = "8"
$ws.Range("D181").Value = "8"

# Row 182 (movie #181) - o baba(ca)
$ws.Range("A151:J151").Copy()
$ws.Range("A182:J182").PasteSpecial(-4122)
$ws.Range("A182").Value = 181
$ws.Range("B182").Value = "o baba(ca)"
$ws.Range("C182").Value = "9"
$ws.Range("D182").Value = "8,5"

# Row 183 (movie #182) - os caça-noivas
$ws.Range("A151:J151").Copy()
$ws.Range("A183:J183").PasteSpecial(-4122)
$ws.Range("A183").Value = 182
$ws.Range("B183").Value = "os caça-noivas"
$ws.Range("C183").Value = "8"
$ws.Range("D183").Value = "8,5"

# Row 184 (movie #183) - a ressaca
$ws.Range("A151:J151").Copy()
$ws.Range("A184:J184").PasteSpecial(-4122)
$ws.Range("A184").Value = 183
$ws.Range("B184").Value = "a ressaca"
$ws.Range("C184").Value = "8"
$ws.Range("D184").Value = "7,5"

# Row 185 (movie #184) - a ressaca 2
$ws.Range("A151:J151").Copy()
$ws.Range("A185:J185").PasteSpecial(-4122)
$ws.Range("A185").Value = 184
$ws.Range("B185").Value = "a ressaca 2"
$ws.Range("C185").Value = "6"
$ws.Range("D185").Value = "pútrido"

# Row 186 (movie #185) - um parto de viagem
$ws.Range("A151:J151").Copy()
$ws.Range("A186:J186").PasteSpecial(-4122)
$ws.Range("A186").Value = 185
$ws.Range("B186").Value = "um parto de viagem"
$ws.Range("C186").Value = "8"
$ws.Range("D186").Value = "7,5"

# Row 187 (movie #186) - o jogo da imitação
$ws.Range("A151:J151").Copy()
$ws.Range("A187:J187").PasteSpecial(-4122)
$ws.Range("A187").Value = 186
$ws.Range("B187").Value = "o jogo da imitação"
$ws.Range("C187").Value = "10"
$ws.Range("D187").Value = "9"

# Row 188 (movie #187) - bem vindo à prisão
$ws.Range("A151:J151").Copy()
$ws.Range("A188:J188").PasteSpecial(-4122)
$ws.Range("A188").Value = 187
$ws.Range("B188").Value = "bem vindo à prisão"
$ws.Range("C188").Value = "8"
$ws.Range("D188").Value = "8,2"

# Row 189 (movie #188) - garota infernal
$ws.Range("A151:J151").Copy()
$ws.Range("A189:J189").PasteSpecial(-4122)
$ws.Range("A189").Value = 188
$ws.Range("B189").Value = "garota infernal"
$ws.Range("C189").Value = "7"
$ws.Range("D189").Value = "7,5"

# Row 190 (movie #189) - defendor
$ws.Range("A151:J151").Copy()
$ws.Range("A190:J190").PasteSpecial(-4122)
$ws.Range("A190").Value = 189
$ws.Range("B190").Value = "defendor"
$ws.Range("C190").Value = "4"
$ws.Range("D190").Value = "3"
$ws.Range("J190").Value = "dropado"

# Row 191 (movie #190) - agente 86
$ws.Range("A151:J151").Copy()
$ws.Range("A191:J191").PasteSpecial(-4122)
$ws.Range("A191").Value = 190
$ws.Range("B191").Value = "agente 86"
$ws.Range("C191").Value = "9"
$ws.Range("D191").Value = "8'"

# Row 192 (movie #191) - Tiras, só que não
$ws.Range("A151:J151").Copy()
$ws.Range("A192:J192").PasteSpecial(-4122)
$ws.Range("A192").Value = 191
$ws.Range("B192").Value = "Tiras, só que não"
$ws.Range("C192").Value = "8,5"
$ws.Range("D192").Value = "8,5"

# Row 193 (movie #192) - homem aranha no aranhaverso
$ws.Range("A151:J151").Copy()
$ws.Range("A193:J193").PasteSpecial(-4122)
$ws.Range("A193").Value = 192
$ws.Range("B193").Value = "homem aranha no aranhaverso"
$ws.Range("C193").Value = "9,8"
$ws.Range("D193").Value = "9"
$ws.Range("F193").Value = "9"

# Row 194 (movie #193) - olhos famintos
$ws.Range("A151:J151").Copy()
$ws.Range("A194:J194").PasteSpecial(-4122)
$ws.Range("A194").Value = 193
$ws.Range("B194").Value = "olhos famintos"
$ws.Range("C194").Value = "8"
$ws.Range("D194").Value = "5"

# Row 195 (movie #194) - O operário
$ws.Range("A151:J151").Copy()
$ws.Range("A195:J195").PasteSpecial(-4122)
$ws.Range("A195").Value = 194
$ws.Range("B195").Value = "O operário"
$ws.Range("C195").Value = "8,5"
$ws.Range("D195").Value = "8,5"

# Row 196 (movie #195) - Halloween
$ws.Range("A151:J151").Copy()
$ws.Range("A196:J196").PasteSpecial(-4122)
$ws.Range("A196").Value = 195
$ws.Range("B196").Value = "Halloween"
$ws.Range("C196").Value = "8"
$ws.Range("D196").Value = "8"

# Row 197 (movie #196) - Homem aranha sem volta pra casa
$ws.Range("A151:J151").Copy()
$ws.Range("A197:J197").PasteSpecial(-4122)
$ws.Range("A197").Value = 196
$ws.Range("B197").Value = "Homem aranha sem volta pra casa"
$ws.Range("C197").Value = "10"
$ws.Range("D197").Value = "11"
$ws.Range("G197").Value = "10"
$ws.Range("J197").Value = "cinema"

$excel.CutCopyMode = $false
